$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

# Row 2
$ws.Range("D2").Value = 0.000249768141657114
$ws.Range("E2").Value = 0.0853085177950561
$ws.Range("G2").Value = 0.00559477461501956
$ws.Range("H2").Value = 0.01019760640338063
$ws.Range("I2").Value = 0.02937948796898127
$ws.Range("J2").Value = 0.03353428887203336
$ws.Range("K2").Value = 0.001744922716170549

# Row 3
$ws.Range("D3").Value = 0.003811206202954054
$ws.Range("E3").Value = 0.06394722079858184
$ws.Range("G3").Value = 0.00335369911044836
$ws.Range("H3").Value = 0.01064082141965628
$ws.Range("I3").Value = 0.01833307696506381
$ws.Range("J3").Value = 0.0274820146150887
$ws.Range("K3").Value = 0.001115125138312578

# Row 4
$ws.Range("D4").Value = 0.003984733484685421
$ws.Range("E4").Value = 0.06672869389876723
$ws.Range("G4").Value = 0.0034317746758461
$ws.Range("H4").Value = 0.01091049751266837
$ws.Range("I4").Value = 0.01965167745947838
$ws.Range("J4").Value = 0.02837291173636913
$ws.Range("K4").Value = 0.001082807779312134

# Row 5
$ws.Range("D5").Value = 0.0003545717336237431
$ws.Range("E5").Value = 0.08475102204829454
$ws.Range("G5").Value = 0.005137935746461153
$ws.Range("H5").Value = 0.009851789567619562
$ws.Range("I5").Value = 0.02968937717378139
$ws.Range("J5").Value = 0.0335434777662158
$ws.Range("K5").Value = 0.001752892509102821

# Row 6
$ws.Range("D6").Value = 0.005934383720159531
$ws.Range("E6").Value = 0.4651480712927878
$ws.Range("G6").Value = 0.01034766295924783
$ws.Range("H6").Value = 0.02766285091638565
$ws.Range("I6").Value = 0.363250554073602
$ws.Range("J6").Value = 0.05025405017659068
$ws.Range("K6").Value = 0.00386454164981842

# Row 8
$ws.Range("D8").Value = 0.000249768141657114
$ws.Range("E8").Value = 0.0853085177950561
$ws.Range("G8").Value = 0.00559477461501956
$ws.Range("H8").Value = 0.01019760640338063
$ws.Range("I8").Value = 0.02937948796898127
$ws.Range("J8").Value = 0.03353428887203336
$ws.Range("K8").Value = 0.001744922716170549

# Row 9
$ws.Range("D9").Value = 0.003811206202954054
$ws.Range("E9").Value = 0.06394722079858184
$ws.Range("G9").Value = 0.00335369911044836
$ws.Range("H9").Value = 0.01064082141965628
$ws.Range("I9").Value = 0.01833307696506381
$ws.Range("J9").Value = 0.0274820146150887
$ws.Range("K9").Value = 0.001115125138312578

# Row 10
$ws.Range("D10").Value = 0.003984733484685421
$ws.Range("E10").Value = 0.06672869389876723
$ws.Range("G10").Value = 0.0034317746758461
$ws.Range("H10").Value = 0.01091049751266837
$ws.Range("I10").Value = 0.01965167745947838
$ws.Range("J10").Value = 0.02837291173636913
$ws.Range("K10").Value = 0.001082807779312134

# Row 11
$ws.Range("D11").Value = 0.0003545717336237431
$ws.Range("E11").Value = 0.08475102204829454
$ws.Range("G11").Value = 0.005137935746461153
$ws.Range("H11").Value = 0.009851789567619562
$ws.Range("I11").Value = 0.02968937717378139
$ws.Range("J11").Value = 0.0335434777662158
$ws.Range("K11").Value = 0.001752892509102821

# Row 12
$ws.Range("D12").Value = 0.005934383720159531
$ws.Range("E12").Value = 0.4651480712927878
$ws.Range("G12").Value = 0.01034766295924783
$ws.Range("H12").Value = 0.02766285091638565
$ws.Range("I12").Value = 0.363250554073602
$ws.Range("J12").Value = 0.05025405017659068
$ws.Range("K12").Value = 0.00386454164981842
